# Auto-generated Excel COM-interop script applying the Chocobo_Profits diff
# across the 8 job-sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 88
$ws.Range("H88").Value2 = 4516.222
$ws.Range("I88").Value2 = 4233
$ws.Range("J88").Value2 = 4572.8667
$ws.Range("K88").Value2 = 4233
$ws.Range("L88").Value2 = 4572.8667
$ws.Range("M88").Value2 = -3827
$ws.Range("N88").Value2 = -5384.8667
# Row 91
$ws.Range("H91").Value2 = 4516.222
$ws.Range("I91").Value2 = 4233
$ws.Range("J91").Value2 = 4572.8667
$ws.Range("K91").Value2 = 4233
$ws.Range("L91").Value2 = 4572.8667
$ws.Range("M91").Value2 = -2829
$ws.Range("N91").Value2 = -7380.8667
# Row 98
$ws.Range("H98").Value2 = 2934.4167
$ws.Range("I98").Value2 = 996.625
$ws.Range("J98").Value2 = 6810
$ws.Range("K98").Value2 = 996.625
$ws.Range("L98").Value2 = 6810
$ws.Range("M98").Value2 = 501.375
$ws.Range("N98").Value2 = -9806
# Row 112
$ws.Range("H112").Value2 = 1330.3877
$ws.Range("I112").Value2 = 649.8570999999999
$ws.Range("J112").Value2 = 1443.8096
$ws.Range("K112").Value2 = 1949.5713
$ws.Range("L112").Value2 = 4331.4288
$ws.Range("M112").Value2 = -841.5712999999998
$ws.Range("N112").Value2 = -6547.4288
# Row 122
$ws.Range("H122").Value2 = 2934.4167
$ws.Range("I122").Value2 = 996.625
$ws.Range("J122").Value2 = 6810
$ws.Range("K122").Value2 = 2989.875
$ws.Range("L122").Value2 = 20430
$ws.Range("M122").Value2 = -539.875
$ws.Range("N122").Value2 = -25330
# Row 132
$ws.Range("H132").Value2 = 575842.75
$ws.Range("I132").Value2 = 326693.22
$ws.Range("K132").Value2 = 980079.6599999999
$ws.Range("M132").Value2 = -977549.6599999999
# Row 135
$ws.Range("H135").Value2 = 308.68967
$ws.Range("I135").Value2 = 214.08
$ws.Range("J135").Value2 = 900
$ws.Range("K135").Value2 = 1926.72
$ws.Range("L135").Value2 = 8100
$ws.Range("M135").Value2 = 608.28
$ws.Range("N135").Value2 = -13170
# Row 138
$ws.Range("H138").Value2 = 1888.2688
$ws.Range("I138").Value2 = 699.63635
$ws.Range("J138").Value2 = 2542.0166
$ws.Range("K138").Value2 = 2098.90905
$ws.Range("L138").Value2 = 7626.0498
$ws.Range("M138").Value2 = 3041.09095
$ws.Range("N138").Value2 = -17906.0498

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value2 = 5963.431
$ws.Range("I32").Value2 = 5357.5776
$ws.Range("K32").Value2 = 5357.5776
$ws.Range("M32").Value2 = -5070.5776
# Row 74
$ws.Range("H74").Value2 = 4891.4546
$ws.Range("I74").Value2 = 4410.095
$ws.Range("J74").Value2 = 15000
$ws.Range("K74").Value2 = 4410.095
$ws.Range("L74").Value2 = 15000
$ws.Range("M74").Value2 = -3536.095
$ws.Range("N74").Value2 = -16748
# Row 77
$ws.Range("H77").Value2 = 4891.4546
$ws.Range("I77").Value2 = 4410.095
$ws.Range("J77").Value2 = 15000
$ws.Range("K77").Value2 = 22050.475
$ws.Range("L77").Value2 = 75000
$ws.Range("M77").Value2 = -17682.475
$ws.Range("N77").Value2 = -83736
# Row 80
$ws.Range("H80").Value2 = 30285.889
$ws.Range("I80").Value2 = 10000
$ws.Range("K80").Value2 = 10000
$ws.Range("M80").Value2 = -9002
# Row 83
$ws.Range("H83").Value2 = 30285.889
$ws.Range("I83").Value2 = 10000
$ws.Range("K83").Value2 = 30000
$ws.Range("M83").Value2 = -25008
# Row 122
$ws.Range("H122").Value2 = 3201
$ws.Range("I122").Value2 = 1803
$ws.Range("J122").Value2 = 4133
$ws.Range("K122").Value2 = 5409
$ws.Range("L122").Value2 = 12399
$ws.Range("M122").Value2 = -2959
$ws.Range("N122").Value2 = -17299
# Row 132
$ws.Range("H132").Value2 = 2300.9722
$ws.Range("I132").Value2 = 1197.6
$ws.Range("J132").Value2 = 4808.636
$ws.Range("K132").Value2 = 3592.8
$ws.Range("L132").Value2 = 14425.908
$ws.Range("M132").Value2 = -1062.8
$ws.Range("N132").Value2 = -19485.908

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 76
$ws.Range("H76").Value2 = 0
$ws.Range("J76").Value2 = 0
$ws.Range("L76").Value2 = 0
$ws.Range("N76").ClearContents()  # was -20630
# Row 79
$ws.Range("H79").Value2 = 0
$ws.Range("J79").Value2 = 0
$ws.Range("L79").Value2 = 0
$ws.Range("N79").ClearContents()  # was -22184
# Row 133
$ws.Range("H133").Value2 = 65000
$ws.Range("J133").Value2 = 100000
$ws.Range("L133").Value2 = 100000
$ws.Range("N133").Value2 = -110120

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value2 = 1867.6438
$ws.Range("I58").Value2 = 1670.9833
$ws.Range("J58").Value2 = 2775.3076
$ws.Range("K58").Value2 = 1670.9833
$ws.Range("L58").Value2 = 2775.3076
$ws.Range("M58").Value2 = -1467.9833
$ws.Range("N58").Value2 = -3181.3076
# Row 107
$ws.Range("H107").Value2 = 669.0833
$ws.Range("I107").Value2 = 573
$ws.Range("J107").Value2 = 1149.5
$ws.Range("K107").Value2 = 573
$ws.Range("L107").Value2 = 1149.5
$ws.Range("M107").Value2 = 1347
$ws.Range("N107").Value2 = -4989.5
# Row 132
$ws.Range("H132").Value2 = 3576.2778
$ws.Range("I132").Value2 = 3236.261
$ws.Range("J132").Value2 = 4177.846
$ws.Range("K132").Value2 = 9708.782999999999
$ws.Range("L132").Value2 = 12533.538
$ws.Range("M132").Value2 = -7178.782999999999
$ws.Range("N132").Value2 = -17593.538
# Row 134
$ws.Range("H134").Value2 = 5757.44
$ws.Range("I134").Value2 = 8534.076999999999
$ws.Range("K134").Value2 = 25602.231
$ws.Range("M134").Value2 = -23067.231
# Row 136
$ws.Range("H136").Value2 = 1867.6438
$ws.Range("I136").Value2 = 1670.9833
$ws.Range("J136").Value2 = 2775.3076
$ws.Range("K136").Value2 = 5012.949900000001
$ws.Range("L136").Value2 = 8325.9228
$ws.Range("M136").Value2 = -2462.949900000001
$ws.Range("N136").Value2 = -13425.9228

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value2 = 103.07143
$ws.Range("I12").Value2 = 23
$ws.Range("J12").Value2 = 163.125
$ws.Range("K12").Value2 = 69
$ws.Range("L12").Value2 = 489.375
$ws.Range("M12").Value2 = 104
$ws.Range("N12").Value2 = -835.375
# Row 131
$ws.Range("H131").Value2 = 8772833
$ws.Range("J131").Value2 = 962.9231
$ws.Range("L131").Value2 = 2888.7693
$ws.Range("N131").Value2 = -12968.7693
# Row 137
$ws.Range("H137").Value2 = 3816.1
$ws.Range("J137").Value2 = 4017.889
$ws.Range("L137").Value2 = 12053.667
$ws.Range("N137").Value2 = -22253.667

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value2 = 50002344
$ws.Range("I80").Value2 = 83335240
$ws.Range("J80").Value2 = 3000
$ws.Range("K80").Value2 = 83335240
$ws.Range("L80").Value2 = 3000
$ws.Range("M80").Value2 = -83334242
$ws.Range("N80").Value2 = -4996
# Row 83
$ws.Range("H83").Value2 = 50002344
$ws.Range("I83").Value2 = 83335240
$ws.Range("J83").Value2 = 3000
$ws.Range("K83").Value2 = 416676200
$ws.Range("L83").Value2 = 15000
$ws.Range("M83").Value2 = -416671208
$ws.Range("N83").Value2 = -24984
# Row 112
$ws.Range("H112").Value2 = 28417.5
$ws.Range("J112").Value2 = 28417.5
$ws.Range("L112").Value2 = 28417.5
$ws.Range("N112").Value2 = -30633.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value2 = 58150
$ws.Range("J36").Value2 = 58150
$ws.Range("L36").Value2 = 58150
$ws.Range("N36").Value2 = -59274
# Row 46
$ws.Range("H46").Value2 = 2023.6428
$ws.Range("J46").Value2 = 1705.2
$ws.Range("L46").Value2 = 1705.2
$ws.Range("N46").Value2 = -2081.2
# Row 133
$ws.Range("H133").Value2 = 0
$ws.Range("J133").Value2 = 0
$ws.Range("L133").Value2 = 0
$ws.Range("N133").ClearContents()  # was -38062.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value2 = 6945826
$ws.Range("I132").Value2 = 883.1905
$ws.Range("K132").Value2 = 2649.5715
$ws.Range("M132").Value2 = -119.5715
# Row 136
$ws.Range("H136").Value2 = 2482.2163
$ws.Range("I136").Value2 = 719.3570999999999
$ws.Range("K136").Value2 = 2158.0713
$ws.Range("M136").Value2 = 391.9287000000004
